# Insert a new data row at row 130 (pushing existing rows 130-214 down to
# 131-215) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("130").Insert()

$ws.Range("A130").Value = 11
$ws.Range("B130").Value = "Vega Monumental Concepción"
$ws.Range("C130").Value = "Bíobío"
$ws.Range("D130").Value = 44904
$ws.Range("E130").Value = 8
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100103
$ws.Range("H130").Value = "Frutos de hueso (carozo)"
$ws.Range("I130").Value = 100103004
$ws.Range("J130").Value = "Durazno"
$ws.Range("K130").Value = "Early Majestic"
$ws.Range("L130").Value = "Primera"
$ws.Range("M130").Value = 100
$ws.Range("N130").Value = 14000
$ws.Range("O130").Value = 15000
$ws.Range("P130").Value = 14500
$ws.Range("Q130").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R130").Value = "Región de O'Higgins"
$ws.Range("S130").Value = 967
$ws.Range("T130").Value = 15
